# chore: update Sheets via scheduled runner
# Refreshes the market-board-derived profit columns (currentAveragePrice*,
# LevePrice*, LeveProfit*) on the per-job leve tables with the latest
# snapshot values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 708
$ws.Range("I6").Value = 435
$ws.Range("J6").Value = 1800
$ws.Range("K6").Value = 1305
$ws.Range("L6").Value = 5400
$ws.Range("M6").Value = -1193
$ws.Range("N6").Value = -5624
$ws.Range("H9").Value = 130.36363
$ws.Range("I9").Value = 123.4
$ws.Range("K9").Value = 123.4
$ws.Range("M9").Value = 45.59999999999999
$ws.Range("H29").Value = 1266.25
$ws.Range("I29").Value = 18.571428
$ws.Range("J29").Value = 10000
$ws.Range("K29").Value = 55.71428400000001
$ws.Range("L29").Value = 30000
$ws.Range("M29").Value = 225.285716
$ws.Range("N29").Value = -30562
$ws.Range("H38").Value = 729.4375
$ws.Range("I38").Value = 67
$ws.Range("K38").Value = 201
$ws.Range("M38").Value = 171
$ws.Range("H98").Value = 10127.546
$ws.Range("I98").Value = 6711.4443
$ws.Range("J98").Value = 25500
$ws.Range("K98").Value = 6711.4443
$ws.Range("L98").Value = 25500
$ws.Range("M98").Value = -5213.4443
$ws.Range("N98").Value = -28496
$ws.Range("H122").Value = 10127.546
$ws.Range("I122").Value = 6711.4443
$ws.Range("J122").Value = 25500
$ws.Range("K122").Value = 20134.3329
$ws.Range("L122").Value = 76500
$ws.Range("M122").Value = -17684.3329
$ws.Range("N122").Value = -81400
$ws.Range("H123").Value = 1729833.4
$ws.Range("J123").Value = 1729833.4
$ws.Range("L123").Value = 1729833.4
$ws.Range("N123").Value = -1739633.4
$ws.Range("H134").Value = 107600
$ws.Range("J134").Value = 107600
$ws.Range("L134").Value = 107600
$ws.Range("N134").Value = -117740
$ws.Range("H136").Value = 77313.336
$ws.Range("J136").Value = 77313.336
$ws.Range("L136").Value = 77313.336
$ws.Range("N136").Value = -87513.336
$ws.Range("H137").Value = 2858.742
$ws.Range("I137").Value = 1856.3
$ws.Range("J137").Value = 3336.0952
$ws.Range("K137").Value = 5568.9
$ws.Range("L137").Value = 10008.2856
$ws.Range("M137").Value = -3018.9
$ws.Range("N137").Value = -15108.2856
$ws.Range("H139").Value = 42545.363
$ws.Range("J139").Value = 42545.363
$ws.Range("L139").Value = 42545.363
$ws.Range("N139").Value = -52825.363
$ws.Range("H140").Value = 122111.11
$ws.Range("J140").Value = 122111.11
$ws.Range("L140").Value = 122111.11
$ws.Range("N140").Value = -132471.11
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11968.441
$ws.Range("I32").Value = 11859.935
$ws.Range("J32").Value = 12914
$ws.Range("K32").Value = 11859.935
$ws.Range("L32").Value = 12914
$ws.Range("M32").Value = -11572.935
$ws.Range("N32").Value = -13488
$ws.Range("H74").Value = 1390.8302
$ws.Range("I74").Value = 1286.05
$ws.Range("J74").Value = 1713.2307
$ws.Range("K74").Value = 1286.05
$ws.Range("L74").Value = 1713.2307
$ws.Range("M74").Value = -412.05
$ws.Range("N74").Value = -3461.2307
$ws.Range("H77").Value = 1390.8302
$ws.Range("I77").Value = 1286.05
$ws.Range("J77").Value = 1713.2307
$ws.Range("K77").Value = 6430.25
$ws.Range("L77").Value = 8566.1535
$ws.Range("M77").Value = -2062.25
$ws.Range("N77").Value = -17302.1535
$ws.Range("H131").Value = 48814.145
$ws.Range("J131").Value = 48814.145
$ws.Range("L131").Value = 48814.145
$ws.Range("N131").Value = -58894.145
$ws.Range("H132").Value = 6040.7715
$ws.Range("I132").Value = 6617.16
$ws.Range("J132").Value = 4599.8
$ws.Range("K132").Value = 19851.48
$ws.Range("L132").Value = 13799.4
$ws.Range("M132").Value = -17321.48
$ws.Range("N132").Value = -18859.4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1960.6118
$ws.Range("I31").Value = 2099.195
$ws.Range("J31").Value = 1831.4773
$ws.Range("K31").Value = 2099.195
$ws.Range("L31").Value = 1831.4773
$ws.Range("M31").Value = -1804.195
$ws.Range("N31").Value = -2421.4773
$ws.Range("H34").Value = 1960.6118
$ws.Range("I34").Value = 2099.195
$ws.Range("J34").Value = 1831.4773
$ws.Range("K34").Value = 2099.195
$ws.Range("L34").Value = 1831.4773
$ws.Range("M34").Value = -1897.195
$ws.Range("N34").Value = -2235.4773
$ws.Range("H58").Value = 2180944.8
$ws.Range("I58").Value = 3369441.8
$ws.Range("J58").Value = 2033.5
$ws.Range("K58").Value = 3369441.8
$ws.Range("L58").Value = 2033.5
$ws.Range("M58").Value = -3369238.8
$ws.Range("N58").Value = -2439.5
$ws.Range("H132").Value = 1129457
$ws.Range("I132").Value = 1691738.1
$ws.Range("J132").Value = 4894.75
$ws.Range("K132").Value = 5075214.300000001
$ws.Range("L132").Value = 14684.25
$ws.Range("M132").Value = -5072684.300000001
$ws.Range("N132").Value = -19744.25
$ws.Range("H134").Value = 1869.8422
$ws.Range("I134").Value = 1565.4706
$ws.Range("K134").Value = 4696.4118
$ws.Range("M134").Value = -2161.4118
$ws.Range("H136").Value = 2180944.8
$ws.Range("I136").Value = 3369441.8
$ws.Range("J136").Value = 2033.5
$ws.Range("K136").Value = 10108325.4
$ws.Range("L136").Value = 6100.5
$ws.Range("M136").Value = -10105775.4
$ws.Range("N136").Value = -11200.5
$ws.Range("H138").Value = 77556
$ws.Range("J138").Value = 77556
$ws.Range("L138").Value = 77556
$ws.Range("N138").Value = -87836
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1268.5051
$ws.Range("I68").Value = 978.5714
$ws.Range("J68").Value = 1564.4791
$ws.Range("K68").Value = 2935.7142
$ws.Range("L68").Value = 4693.4373
$ws.Range("M68").Value = -2124.7142
$ws.Range("N68").Value = -6315.4373
$ws.Range("H71").Value = 1268.5051
$ws.Range("I71").Value = 978.5714
$ws.Range("J71").Value = 1564.4791
$ws.Range("K71").Value = 8807.142600000001
$ws.Range("L71").Value = 14080.3119
$ws.Range("M71").Value = -4751.142600000001
$ws.Range("N71").Value = -22192.3119
$ws.Range("H107").Value = 1280
$ws.Range("I107").Value = 1065.2653
$ws.Range("J107").Value = 1758.2727
$ws.Range("K107").Value = 3195.7959
$ws.Range("L107").Value = 5274.8181
$ws.Range("M107").Value = -1275.7959
$ws.Range("N107").Value = -9114.8181
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 10951.667
$ws.Range("J109").Value = 10951.667
$ws.Range("L109").Value = 10951.667
$ws.Range("N109").Value = -13031.667
$ws.Range("H132").Value = 3293.1765
$ws.Range("I132").Value = 2314
$ws.Range("J132").Value = 4163.5557
$ws.Range("K132").Value = 6942
$ws.Range("L132").Value = 12490.6671
$ws.Range("M132").Value = -4412
$ws.Range("N132").Value = -17550.6671
$ws.Range("H140").Value = 60225.715
$ws.Range("J140").Value = 60225.715
$ws.Range("L140").Value = 60225.715
$ws.Range("N140").Value = -70585.715
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4088
$ws.Range("I136").Value = 2102
$ws.Range("J136").Value = 4750
$ws.Range("K136").Value = 6306
$ws.Range("L136").Value = 14250
$ws.Range("N136").Value = -19350
$ws.Range("M136").Value = -3756
$ws.Range("H137").Value = 37428.625
$ws.Range("J137").Value = 37428.625
$ws.Range("L137").Value = 37428.625
$ws.Range("N137").Value = -47628.625
$ws.Range("H140").Value = 59582.25
$ws.Range("J140").Value = 59582.25
$ws.Range("L140").Value = 59582.25
$ws.Range("N140").Value = -69942.25
$ws.Range("H141").Value = 49598
$ws.Range("J141").Value = 49598
$ws.Range("L141").Value = 49598
$ws.Range("N141").Value = -59958
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 15627208
$ws.Range("I122").Value = 15627208
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 46881624
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -46879174
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 3199.1428
$ws.Range("I132").Value = 2977.875
$ws.Range("K132").Value = 8933.625
$ws.Range("M132").Value = -6403.625
$ws.Range("H136").Value = 3466.4583
$ws.Range("I136").Value = 3766.25
$ws.Range("J136").Value = 3166.6667
$ws.Range("K136").Value = 11298.75
$ws.Range("L136").Value = 9500.000100000001
$ws.Range("M136").Value = -8748.75
$ws.Range("N136").Value = -14600.0001
$ws.Range("H139").Value = 47926.875
$ws.Range("J139").Value = 47630.715
$ws.Range("L139").Value = 47630.715
$ws.Range("N139").Value = -57910.715
